$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feed")

# Thread pitch / feed-rate table edits (column A = value, column D = unit "mm"/"tpi")
$ws.Range("C6").Value = 1.5

$ws.Range("A9").Value = 0.09

$ws.Range("A10").Value = 0.15
$ws.Range("D10").Value = "mm"

$ws.Range("A11").Value = 1.5

$ws.Range("A12").Value = 1

$ws.Range("A13").Value = 10
$ws.Range("D13").Value = "tpi"

# Update the active selection to reflect the cell the author last edited
$ws.Range("I17").Select()
